# Adatum Leave of Absence Policy - split 4 runs around a proofing (grammar)
# mark, matching the author's edit. Each target run of plain text is split
# into three runs: "<lead> ", "<flagged word(s)>", " <tail>" with a
# <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/> pair
# bracketing the middle run. Paragraph identity (w14:paraId/w:rsidR/etc.)
# and any <w:pPr> block are preserved exactly as in the source paragraph.

$d = $word.ActiveDocument

function Get-ParaPrefix {
    param([string]$Xml, [int]$Idx)

    $sub = $Xml.Substring(0, $Idx)
    $lastOpen = $sub.LastIndexOf("<w:p>")
    $lastOpen2 = $sub.LastIndexOf("<w:p ")
    if ($lastOpen2 -gt $lastOpen) { $lastOpen = $lastOpen2 }

    $rIdx = $Xml.IndexOf("<w:r>", $lastOpen)
    $rIdx2 = $Xml.IndexOf("<w:r ", $lastOpen)
    if ($rIdx -eq -1) { $rIdx = $rIdx2 }
    elseif ($rIdx2 -ne -1 -and $rIdx2 -lt $rIdx) { $rIdx = $rIdx2 }

    return $Xml.Substring($lastOpen, $rIdx - $lastOpen)
}

function Split-WithGrammarMark {
    param(
        [string]$Whole,
        [string]$Lead,
        [string]$Flag,
        [string]$Tail,
        [string]$RprXml
    )

    # Locate the paragraph's opening tag (+ <w:pPr> if any) so it can be
    # reproduced verbatim in the replacement fragment.
    $fullXml = $d.WordOpenXML
    $needleIdx = $fullXml.IndexOf($Whole)
    if ($needleIdx -lt 0) {
        throw "Could not locate paragraph prefix for: $Whole"
    }
    $prefix = Get-ParaPrefix $fullXml $needleIdx

    # Locate the precise character range of the target text.
    $rng = $d.Content
    $found = $rng.Find.Execute($Whole, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find target text: $Whole"
    }
    $target = $d.Range($rng.Start, $rng.End)

    $r1 = "<w:r>" + $RprXml + "<w:t xml:space=`"preserve`">" + $Lead + "</w:t></w:r>"
    $r2 = "<w:r>" + $RprXml + "<w:t>" + $Flag + "</w:t></w:r>"
    $r3 = "<w:r>" + $RprXml + "<w:t xml:space=`"preserve`">" + $Tail + "</w:t></w:r>"

    $body = $r1 + '<w:proofErr w:type="gramStart"/>' + $r2 + '<w:proofErr w:type="gramEnd"/>' + $r3

    $frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
            '<w:body>' + $prefix + $body + '</w:p></w:body></w:document>' +
            '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($frag)
}

$italicRpr = "<w:rPr><w:i/></w:rPr>"

Split-WithGrammarMark "Effective: October 31, 2025 | v1.0" "Effective: October 31, " "2025" " | v1.0" $italicRpr

Split-WithGrammarMark "Job-protected where eligible; may coordinate with STD." "Job-protected where " "eligible;" " may coordinate with STD." ""

Split-WithGrammarMark "Care for qualifying family member with serious health condition." "Care for qualifying family " "member" " with serious health condition." ""

Split-WithGrammarMark "Q: What if I exhaust PTO? A: Discuss unpaid time options with your manager and HR." "Q: What if I " "exhaust" " PTO? A: Discuss unpaid time options with your manager and HR." ""

Write-Output "edits applied"
